$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reassign the table contents in place (no column shifting needed - the final
# layout still only spans columns A:D, just with different headers/values).

# Row 1 headers
$ws.Range("A1").Value = "Branch_Num"
$ws.Range("B1").Value = "Location"
$ws.Range("C1").Value = "Num_Employees"
$ws.Range("D1").Value = "Num_Clients"

# Row 2 data
$ws.Range("A2").Value = 20240519
$ws.Range("B2").Value = "Victoria"
$ws.Range("C2").Value = 15
$ws.Range("D2").Value = 100

# Column widths - columns B and D already have the right width and are left
# untouched; only the new column A and the changed column C need updating.
$ws.Columns.Item(1).ColumnWidth = 14.42
$ws.Columns.Item(3).ColumnWidth = 16.92

# Restore active selection to match saved state
$ws.Range("E7").Select()
